$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (B, C) ---
$ws.Columns.Item(2).ColumnWidth = 35.451822916666664
$ws.Columns.Item(3).ColumnWidth = 26.592447916666668

# --- Values ---
$ws.Range("B2").Value = "TAAK"
$ws.Range("C2").Value = "UITGEVOERD DOOR"

$ws.Range("B3").Value = "Base repo aanmaken en opvullen"
$ws.Range("C3").Value = "Thomas"

$ws.Range("B4").Value = "DNS script providen en aanpassen"
$ws.Range("C4").Value = "Jelle"

$ws.Range("B5").Value = "Changes aan base scripts"
$ws.Range("C5").Value = "Thomas "

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 20.25
$ws.Rows.Item(3).RowHeight = 15.75

# --- Normalize alignment flag across the data body BEFORE border work so the
#     border pass lands on the final cell style directly (avoids leaving
#     behind orphan intermediate styles). ---
$ws.Range("B3:C37").WrapText = $false

# --- Top border separating header from first data row (thick, accent1 color) ---
$topBorder = $ws.Range("B3:C3").Borders.Item(8)
$topBorder.Weight = 4
$topBorder.Color = 13998939

# --- Right border separating column B from column C for all data rows ---
$rightBorder = $ws.Range("B3:B37").Borders.Item(10)
$rightBorder.Weight = 2

# --- Header style (B2:C2) : bold "Heading 1" cell style ---
$ws.Range("B2:C2").Style = "Heading 1"

# --- Selection ---
$ws.Range("C9").Select()
